# Generate Report for Handoff
#
# Localization status report: refresh the "Latest Handoff Datetime" (column D)
# for the row that was just handed off (99319a3b-4a97-4338-b29b-ed4a462b98fc.md,
# row 5) on both the zh-cn and de-de target-language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-26 06:07:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-26 06:07:37"
